# Update the EV simulation table on sheet "a" (columns I..U, rows 1..7)
# with the new values, and set columns V..Z (the out-of-range simulation
# years) to the #N/A error value, matching the refreshed model output.
# Also refresh the derived summary values in C10:D15 which are pasted
# values (not formulas) mirroring columns I and the Hoja2 helper
# calculations.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("a")

# New values for columns I (9) through U (21) for rows 1-7.
$rowValues = @{
    1 = @(155, 304, 453, 582, 721.5, 845.5, 974, 1081.5, 1186, 1298, 1409.5, 1515.5, 1602)
    2 = @(155, 304, 453, 582, 721.5, 845.5, 974, 1081.5, 1186, 1298, 1409.5, 1515.5, 1602)
    3 = @(148.5, 284.5, 405, 536, 661, 778, 892, 1010.5, 1122.5, 1215, 1295, 1389, 1480.5)
    4 = @(148.5, 284.5, 405, 536, 661, 778, 892, 1010.5, 1122.5, 1215, 1295, 1389, 1480.5)
    5 = @(101, 208, 307, 403, 491.5, 586.5, 657, 739, 824.5, 887.5, 955.5, 1009.5, 1071)
    6 = @(101, 208, 307, 403, 491.5, 586.5, 657, 739, 824.5, 887.5, 955.5, 1009.5, 1071)
    7 = @(429, 836, 1243, 1629, 1976, 2313, 2653, 2961, 3240, 3548, 3834, 4093, 4345)
}

foreach ($r in $rowValues.Keys) {
    $vals = $rowValues[$r]
    $col = 9   # column I
    foreach ($v in $vals) {
        $ws.Cells.Item($r, $col).Value = $v
        $col = $col + 1
    }
    # Columns V..Z (22..26) now fall outside the simulation horizon -> #N/A
    for ($col = 22; $col -le 26; $col++) {
        $ws.Cells.Item($r, $col).Value = "#N/A"
    }
}

# Refresh the pasted summary values in C10:D15 (these mirror columns
# I/J on sheet "a" multiplied out, but are stored as plain values,
# not formulas, so they must be updated explicitly).
$ws.Range("C10").Value = 7595
$ws.Range("D10").Value = 62

$ws.Range("C11").Value = 4976.4939839999997
$ws.Range("D11").Value = 62

$ws.Range("C14").Value = 11586.280616064001
$ws.Range("D14").Value = 97

$ws.Range("C15").Value = 11009
$ws.Range("D15").Value = 88
